# weight_lift/20241223094620.xlsx — "Add files via upload"
#
# The workbook was re-saved by a different Excel build: the legacy text
# query-table import was converted to plain values, mojibake in the
# sensor-unit headers (and the firmware version string) was corrected,
# and the worksheet was renamed to match the imported file name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

# --- Fix mis-decoded degree-sign ("¬∞" -> "°") in the gyro / angle / ---
# --- temperature column headers living in the shared-string table.  ---
$rng.Replace("AsX(¬∞/s)", "AsX(°/s)")
$rng.Replace("AsY(¬∞/s)", "AsY(°/s)")
$rng.Replace("AsZ(¬∞/s)", "AsZ(°/s)")
$rng.Replace("AngleX(¬∞)", "AngleX(°)")
$rng.Replace("AngleY(¬∞)", "AngleY(°)")
$rng.Replace("AngleZ(¬∞)", "AngleZ(°)")
$rng.Replace("Temperature(¬∞C)", "Temperature(°C)")

# --- Firmware/version string used a comma-separated decimal; normalise ---
# --- it to the dotted form ("10080,1,13" -> "10080.1.13").             ---
$rng.Replace("10080,1,13", "10080.1.13")

# --- The sheet was created from an import named after the source log  ---
# --- file; rename it to match (was the generic "Лист1").              ---
$ws.Name = "20241223094620"

# --- Drop the defined name that pointed at the old query-table range; ---
# --- the data is now a plain range, not an external-data import.      ---
$names = @($wb.Names)
foreach ($n in $names) {
  $n.Delete()
}
